# Auto-generated edit script: updates market-price-derived columns (H-N)
# on the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to reflect a refreshed
# market data pull ("chore: update Sheets via scheduled runner").

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 724.2
$ws.Range("I6").Value = 386.3
$ws.Range("K6").Value = 1158.9
$ws.Range("M6").Value = -1046.9
$ws.Range("H33").Value = 107.56
$ws.Range("I33").Value = 114.55556
$ws.Range("J33").Value = 89.57143000000001
$ws.Range("K33").Value = 114.55556
$ws.Range("L33").Value = 89.57143000000001
$ws.Range("M33").Value = 114.44444
$ws.Range("N33").Value = -547.57143
$ws.Range("H38").Value = 241.75
$ws.Range("J38").Value = 421.16666
$ws.Range("L38").Value = 1263.49998
$ws.Range("N38").Value = -2007.49998
$ws.Range("H40").Value = 2248.25
$ws.Range("I40").Value = 2374.2856
$ws.Range("K40").Value = 2374.2856
$ws.Range("M40").Value = -2199.2856
$ws.Range("H42").Value = 242.8
$ws.Range("I42").Value = 113
$ws.Range("J42").Value = 329.33334
$ws.Range("K42").Value = 339
$ws.Range("L42").Value = 988.0000200000001
$ws.Range("M42").Value = -109
$ws.Range("N42").Value = -1448.00002
$ws.Range("H51").Value = 4412.5
$ws.Range("I51").Value = 3200
$ws.Range("J51").Value = 4692.3076
$ws.Range("K51").Value = 3200
$ws.Range("L51").Value = 4692.3076
$ws.Range("M51").Value = -2716
$ws.Range("N51").Value = -5660.3076
$ws.Range("H64").Value = 4664.636
$ws.Range("J64").Value = 5305.857
$ws.Range("L64").Value = 5305.857
$ws.Range("N64").Value = -5801.857
$ws.Range("H67").Value = 4664.636
$ws.Range("J67").Value = 5305.857
$ws.Range("L67").Value = 5305.857
$ws.Range("N67").Value = -7021.857
$ws.Range("H101").Value = 506
$ws.Range("I101").Value = 487.14285
$ws.Range("J101").Value = 522.5
$ws.Range("K101").Value = 1461.42855
$ws.Range("L101").Value = 1567.5
$ws.Range("M101").Value = 160.5714499999999
$ws.Range("N101").Value = -4811.5
$ws.Range("H115").Value = 1809.4445
$ws.Range("I115").Value = 783.5714
$ws.Range("J115").Value = 5400
$ws.Range("K115").Value = 2350.7142
$ws.Range("L115").Value = 16200
$ws.Range("M115").Value = -783.7142000000003
$ws.Range("N115").Value = -19334
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H138").Value = 8873661
$ws.Range("I138").Value = 2781012
$ws.Range("J138").Value = 11366108
$ws.Range("K138").Value = 8343036
$ws.Range("L138").Value = 34098324
$ws.Range("M138").Value = -8337896
$ws.Range("N138").Value = -34108604
$ws.Range("H141").Value = 3749.2144
$ws.Range("I141").Value = 2616.2727
$ws.Range("K141").Value = 7848.8181
$ws.Range("M141").Value = -2668.8181

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3132.2563
$ws.Range("I61").Value = 2627.25
$ws.Range("J61").Value = 5440.857
$ws.Range("K61").Value = 2627.25
$ws.Range("L61").Value = 5440.857
$ws.Range("M61").Value = -2415.25
$ws.Range("N61").Value = -5864.857
$ws.Range("H103").Value = 35362
$ws.Range("J103").Value = 35362
$ws.Range("L103").Value = 35362
$ws.Range("N103").Value = -37706
$ws.Range("H136").Value = 3132.2563
$ws.Range("I136").Value = 2627.25
$ws.Range("J136").Value = 5440.857
$ws.Range("K136").Value = 7881.75
$ws.Range("L136").Value = 16322.571
$ws.Range("M136").Value = -5331.75
$ws.Range("N136").Value = -21422.571

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H134").Value = 3803.4075
$ws.Range("I134").Value = 2748.2666
$ws.Range("J134").Value = 5122.3335
$ws.Range("K134").Value = 8244.799800000001
$ws.Range("L134").Value = 15367.0005
$ws.Range("M134").Value = -5709.799800000001
$ws.Range("N134").Value = -20437.0005

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5397.5947
$ws.Range("I31").Value = 1970.0869
$ws.Range("J31").Value = 11028.5
$ws.Range("K31").Value = 1970.0869
$ws.Range("L31").Value = 11028.5
$ws.Range("M31").Value = -1675.0869
$ws.Range("N31").Value = -11618.5
$ws.Range("H34").Value = 5397.5947
$ws.Range("I34").Value = 1970.0869
$ws.Range("J34").Value = 11028.5
$ws.Range("K34").Value = 1970.0869
$ws.Range("L34").Value = 11028.5
$ws.Range("M34").Value = -1768.0869
$ws.Range("N34").Value = -11432.5
$ws.Range("H58").Value = 2256.7727
$ws.Range("I58").Value = 1636.6666
$ws.Range("J58").Value = 3585.5715
$ws.Range("K58").Value = 1636.6666
$ws.Range("L58").Value = 3585.5715
$ws.Range("M58").Value = -1433.6666
$ws.Range("N58").Value = -3991.5715
$ws.Range("H70").Value = 29000
$ws.Range("J70").Value = 29000
$ws.Range("L70").Value = 29000
$ws.Range("N70").Value = -29630
$ws.Range("H73").Value = 29000
$ws.Range("J73").Value = 29000
$ws.Range("L73").Value = 29000
$ws.Range("N73").Value = -31184
$ws.Range("H106").Value = 30000
$ws.Range("J106").Value = 30000
$ws.Range("L106").Value = 30000
$ws.Range("N106").Value = -32524
$ws.Range("H132").Value = 2718.5625
$ws.Range("I132").Value = 2466.4443
$ws.Range("J132").Value = 3474.9167
$ws.Range("K132").Value = 7399.3329
$ws.Range("L132").Value = 10424.7501
$ws.Range("M132").Value = -4869.3329
$ws.Range("N132").Value = -15484.7501
$ws.Range("H136").Value = 2256.7727
$ws.Range("I136").Value = 1636.6666
$ws.Range("J136").Value = 3585.5715
$ws.Range("K136").Value = 4909.9998
$ws.Range("L136").Value = 10756.7145
$ws.Range("M136").Value = -2359.9998
$ws.Range("N136").Value = -15856.7145

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 88.25
$ws.Range("I7").Value = 88.25
$ws.Range("K7").Value = 264.75
$ws.Range("M7").Value = -152.75
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H80").Value = 1198.25
$ws.Range("I80").Value = 390
$ws.Range("J80").Value = 1467.6666
$ws.Range("K80").Value = 1170
$ws.Range("L80").Value = 4402.9998
$ws.Range("M80").Value = -234
$ws.Range("N80").Value = -6274.9998
$ws.Range("H83").Value = 1198.25
$ws.Range("I83").Value = 390
$ws.Range("J83").Value = 1467.6666
$ws.Range("K83").Value = 3510
$ws.Range("L83").Value = 13208.9994
$ws.Range("M83").Value = 1170
$ws.Range("N83").Value = -22568.9994
$ws.Range("H92").Value = 1000
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 3000
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -5496
$ws.Range("H122").Value = 1202.1111
$ws.Range("I122").Value = 256.22223
$ws.Range("J122").Value = 2148
$ws.Range("K122").Value = 2306.00007
$ws.Range("L122").Value = 19332
$ws.Range("M122").Value = 143.9999299999999
$ws.Range("N122").Value = -24232

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2046.9714
$ws.Range("I132").Value = 1842.1428
$ws.Range("J132").Value = 2866.2856
$ws.Range("K132").Value = 5526.428400000001
$ws.Range("L132").Value = 8598.856800000001
$ws.Range("M132").Value = -2996.428400000001
$ws.Range("N132").Value = -13658.8568

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2788.75
$ws.Range("I40").Value = 1571.6666
$ws.Range("J40").Value = 3194.4443
$ws.Range("K40").Value = 1571.6666
$ws.Range("L40").Value = 3194.4443
$ws.Range("M40").Value = -1435.6666
$ws.Range("N40").Value = -3466.4443
$ws.Range("H55").Value = 488.5
$ws.Range("I55").Value = 619.5
$ws.Range("J55").Value = 455.75
$ws.Range("K55").Value = 619.5
$ws.Range("L55").Value = 455.75
$ws.Range("M55").Value = -446.5
$ws.Range("N55").Value = -801.75
$ws.Range("H96").Value = 34000
$ws.Range("J96").Value = 34000
$ws.Range("L96").Value = 34000
$ws.Range("N96").Value = -39492

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2719.7058
$ws.Range("I132").Value = 2372.5652
$ws.Range("J132").Value = 3445.5454
$ws.Range("K132").Value = 7117.6956
$ws.Range("L132").Value = 10336.6362
$ws.Range("M132").Value = -4587.6956
$ws.Range("N132").Value = -15396.6362
